$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column L (last existing year column) to new column M
$ws.Range("L3:L5").Copy()
$ws.Range("M3:M5").PasteSpecial(-4122)  # xlPasteFormats

# Add new column M data for year 2022
$ws.Range("M4").Value = 2022
$ws.Range("M5").Value = 373

# Update selection to match target
$ws.Range("O4").Select()
